$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("K12").Value = 149.5
$ws.Range("H12").Value = 149.5
$ws.Range("I12").Value = 149.5
$ws.Range("M12").Value = 20.5
$ws.Range("M86").Value = -8381.5
$ws.Range("K86").Value = 9504.5
$ws.Range("H86").Value = 9504.5
$ws.Range("I86").Value = 9504.5
$ws.Range("M88").Value = -2111.6667
$ws.Range("N88").Value = -6412.8
$ws.Range("K88").Value = 2517.6667
$ws.Range("H88").Value = 4444.625
$ws.Range("I88").Value = 2517.6667
$ws.Range("L88").Value = 5600.8
$ws.Range("J88").Value = 5600.8
$ws.Range("K89").Value = 47522.5
$ws.Range("H89").Value = 9504.5
$ws.Range("I89").Value = 9504.5
$ws.Range("M89").Value = -41906.5
$ws.Range("L91").Value = 5600.8
$ws.Range("J91").Value = 5600.8
$ws.Range("H91").Value = 4444.625
$ws.Range("I91").Value = 2517.6667
$ws.Range("M91").Value = -1113.6667
$ws.Range("N91").Value = -8408.799999999999
$ws.Range("K91").Value = 2517.6667
$ws.Range("L98").Value = 11625
$ws.Range("J98").Value = 11625
$ws.Range("H98").Value = 6553.0415
$ws.Range("I98").Value = 1481.0834
$ws.Range("M98").Value = 16.91660000000002
$ws.Range("N98").Value = -14621
$ws.Range("K98").Value = 1481.0834
$ws.Range("K122").Value = 4443.2502
$ws.Range("J122").Value = 11625
$ws.Range("H122").Value = 6553.0415
$ws.Range("I122").Value = 1481.0834
$ws.Range("L122").Value = 34875
$ws.Range("N122").Value = -39775
$ws.Range("M122").Value = -1993.2502
$ws.Range("M138").Value = 1567
$ws.Range("K138").Value = 3573
$ws.Range("H138").Value = 2224.7144
$ws.Range("I138").Value = 1191
$ws = $wb.Worksheets.Item(2)
$ws.Range("M32").Value = -163.5
$ws.Range("K32").Value = 450.5
$ws.Range("H32").Value = 450.5
$ws.Range("I32").Value = 450.5
$ws.Range("H63").Value = 3301.1
$ws.Range("I63").Value = 3446.2222
$ws.Range("M63").Value = -2760.2222
$ws.Range("N63").Value = -3367
$ws.Range("K63").Value = 3446.2222
$ws.Range("L63").Value = 1995
$ws.Range("J63").Value = 1995
$ws.Range("K66").Value = 17231.111
$ws.Range("N66").Value = -16839
$ws.Range("L66").Value = 9975
$ws.Range("H66").Value = 3301.1
$ws.Range("I66").Value = 3446.2222
$ws.Range("M66").Value = -13799.111
$ws.Range("H88").Value = 2540.4
$ws.Range("H91").Value = 2540.4
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("J125").Value = 0
$ws.Range("H125").Value = 0
$ws = $wb.Worksheets.Item(3)
$ws.Range("M23").Value = 258
$ws.Range("K23").Value = 25
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 25
$ws = $wb.Worksheets.Item(4)
$ws.Range("L2").Value = 495
$ws.Range("J2").Value = 495
$ws.Range("H2").Value = 218.6
$ws.Range("I2").Value = 34.333332
$ws.Range("M2").Value = 78.666668
$ws.Range("N2").Value = -721
$ws.Range("K2").Value = 34.333332
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("M11").Value = 40
$ws.Range("K11").Value = 100
$ws.Range("H22").Value = 766.13336
$ws.Range("I22").Value = 749.2857
$ws.Range("M22").Value = -399.2857
$ws.Range("K22").Value = 749.2857
$ws.Range("H58").Value = 9872.286
$ws.Range("N58").Value = -10713.692
$ws.Range("L58").Value = 10307.692
$ws.Range("J58").Value = 10307.692
$ws.Range("J104").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("K122").Value = 4709.142599999999
$ws.Range("J122").Value = 0
$ws.Range("H122").Value = 1569.7142
$ws.Range("I122").Value = 1569.7142
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("M122").Value = -2259.142599999999
$ws.Range("J136").Value = 10307.692
$ws.Range("H136").Value = 9872.286
$ws.Range("L136").Value = 30923.076
$ws.Range("N136").Value = -36023.076
$ws = $wb.Worksheets.Item(5)
$ws.Range("H50").Value = 387.5
$ws.Range("I50").Value = 300
$ws.Range("M50").Value = -419
$ws.Range("K50").Value = 900
$ws.Range("H53").Value = 387.5
$ws.Range("I53").Value = 300
$ws.Range("M53").Value = -419
$ws.Range("K53").Value = 900
$ws.Range("I121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("K121").Value = 0
$ws.Range("H121").Value = 1603
$ws.Range("H132").Value = 1352.6
$ws.Range("I132").Value = 1287.8
$ws.Range("M132").Value = -9060.199999999999
$ws.Range("K132").Value = 11590.2
$ws = $wb.Worksheets.Item(6)
$ws.Range("N109").Value = -36722.5
$ws.Range("J109").Value = 34642.5
$ws.Range("H109").Value = 34642.5
$ws.Range("L109").Value = 34642.5
$ws.Range("M113").Value = -11579
$ws.Range("K113").Value = 13749
$ws.Range("H113").Value = 10624.75
$ws.Range("I113").Value = 13749
$ws.Range("H132").Value = 7155.625
$ws.Range("I132").Value = 5579.5
$ws.Range("L132").Value = 29347.5
$ws.Range("J132").Value = 9782.5
$ws.Range("M132").Value = -14208.5
$ws.Range("N132").Value = -34407.5
$ws.Range("K132").Value = 16738.5
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 9180.546
$ws.Range("I22").Value = 9000
$ws.Range("M22").Value = -8705
$ws.Range("K22").Value = 9000
$ws.Range("M27").Value = -8893
$ws.Range("K27").Value = 9000
$ws.Range("H27").Value = 9180.546
$ws.Range("I27").Value = 9000
$ws.Range("L140").Value = 89950
$ws.Range("J140").Value = 89950
$ws.Range("H140").Value = 84973.5
$ws.Range("N140").Value = -100310
$ws = $wb.Worksheets.Item(8)
$ws.Range("L14").Value = 3249.9
$ws.Range("J14").Value = 3249.9
$ws.Range("H14").Value = 3249.9
$ws.Range("N14").Value = -3585.9
$ws.Range("N55").Value = -9442
$ws.Range("J55").Value = 8888
$ws.Range("H55").Value = 8888
$ws.Range("L55").Value = 8888
$ws.Range("K122").Value = 9642
$ws.Range("H122").Value = 3214
$ws.Range("I122").Value = 3214
$ws.Range("M122").Value = -7192
$ws.Range("H132").Value = 10228.357
$ws.Range("I132").Value = 6839.6
$ws.Range("M132").Value = -17988.8
$ws.Range("K132").Value = 20518.8
